$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the self-bearing diagonal (buoy to itself) with 0.
# Row 2 -> column B (index 2), row 3 -> column C (index 3), ... row 17 -> column Q (index 17)
for ($row = 2; $row -le 17; $row++) {
    $col = $row
    $ws.Cells.Item($row, $col).Value = 0
}

# Update the active selection left by the editor.
$ws.Range("I19").Select()
